$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1708.1342
$ws.Range("I15").Value = 1708.1342
$ws.Range("K15").Value = 5124.402599999999
$ws.Range("M15").Value = -4955.402599999999
$ws.Range("H32").Value = 1271.5
$ws.Range("I32").Value = 380
$ws.Range("J32").Value = 1398.8572
$ws.Range("K32").Value = 380
$ws.Range("L32").Value = 1398.8572
$ws.Range("M32").Value = -54
$ws.Range("N32").Value = -2050.8572
$ws.Range("H70").Value = 1672.439
$ws.Range("I70").Value = 1655.8518
$ws.Range("J70").Value = 1704.4286
$ws.Range("K70").Value = 4967.555399999999
$ws.Range("L70").Value = 5113.2858
$ws.Range("M70").Value = -4697.555399999999
$ws.Range("N70").Value = -5653.2858
$ws.Range("H73").Value = 1672.439
$ws.Range("I73").Value = 1655.8518
$ws.Range("J73").Value = 1704.4286
$ws.Range("K73").Value = 4967.555399999999
$ws.Range("L73").Value = 5113.2858
$ws.Range("M73").Value = -4031.555399999999
$ws.Range("N73").Value = -6985.2858
$ws.Range("H76").Value = 3553.2646
$ws.Range("I76").Value = 3276.1035
$ws.Range("J76").Value = 5160.8
$ws.Range("K76").Value = 3276.1035
$ws.Range("L76").Value = 5160.8
$ws.Range("M76").Value = -2961.1035
$ws.Range("N76").Value = -5790.8
$ws.Range("H79").Value = 3553.2646
$ws.Range("I79").Value = 3276.1035
$ws.Range("J79").Value = 5160.8
$ws.Range("K79").Value = 3276.1035
$ws.Range("L79").Value = 5160.8
$ws.Range("M79").Value = -2184.1035
$ws.Range("N79").Value = -7344.8
$ws.Range("H132").Value = 2197.2258
$ws.Range("I132").Value = 1637.5555
$ws.Range("K132").Value = 4912.666499999999
$ws.Range("M132").Value = -2382.666499999999
$ws.Range("H141").Value = 7125
$ws.Range("I141").Value = 7166.6665
$ws.Range("J141").Value = 7000
$ws.Range("K141").Value = 21499.9995
$ws.Range("L141").Value = 21000
$ws.Range("M141").Value = -16319.9995
$ws.Range("N141").Value = -31360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3256.3262
$ws.Range("I132").Value = 2939.24
$ws.Range("K132").Value = 8817.719999999999
$ws.Range("M132").Value = -6287.719999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 5790.1665
$ws.Range("I54").Value = 3281.3333
$ws.Range("J54").Value = 13316.667
$ws.Range("K54").Value = 3281.3333
$ws.Range("L54").Value = 13316.667
$ws.Range("M54").Value = -2797.3333
$ws.Range("N54").Value = -14284.667
$ws.Range("H134").Value = 2608.75
$ws.Range("I134").Value = 2547.4546
$ws.Range("J134").Value = 2648.4119
$ws.Range("K134").Value = 7642.3638
$ws.Range("L134").Value = 7945.2357
$ws.Range("M134").Value = -5107.3638
$ws.Range("N134").Value = -13015.2357
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3003.7551
$ws.Range("I31").Value = 1512.1052
$ws.Range("J31").Value = 3948.4666
$ws.Range("K31").Value = 1512.1052
$ws.Range("L31").Value = 3948.4666
$ws.Range("M31").Value = -1217.1052
$ws.Range("N31").Value = -4538.4666
$ws.Range("H34").Value = 3003.7551
$ws.Range("I34").Value = 1512.1052
$ws.Range("J34").Value = 3948.4666
$ws.Range("K34").Value = 1512.1052
$ws.Range("L34").Value = 3948.4666
$ws.Range("M34").Value = -1310.1052
$ws.Range("N34").Value = -4352.4666
$ws.Range("H132").Value = 2801.3845
$ws.Range("I132").Value = 2268.889
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 6806.667
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -4276.667
$ws.Range("N132").Value = -17058.5
$ws.Range("H134").Value = 3584.1538
$ws.Range("I134").Value = 3810.5264
$ws.Range("J134").Value = 2969.7144
$ws.Range("K134").Value = 11431.5792
$ws.Range("L134").Value = 8909.143199999999
$ws.Range("M134").Value = -8896.5792
$ws.Range("N134").Value = -13979.1432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 13471.818
$ws.Range("I3").Value = 18041.428
$ws.Range("J3").Value = 5475
$ws.Range("K3").Value = 54124.284
$ws.Range("L3").Value = 16425
$ws.Range("M3").Value = -54012.284
$ws.Range("N3").Value = -16649
$ws.Range("H107").Value = 689.7308
$ws.Range("I107").Value = 217.33333
$ws.Range("J107").Value = 1094.6428
$ws.Range("K107").Value = 651.99999
$ws.Range("L107").Value = 3283.9284
$ws.Range("M107").Value = 1268.00001
$ws.Range("N107").Value = -7123.928400000001
$ws.Range("H122").Value = 385.57144
$ws.Range("I122").Value = 338.30768
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 3044.76912
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -594.7691199999999
$ws.Range("N122").Value = -13900
$ws.Range("H133").Value = 46366.58
$ws.Range("I133").Value = 252007.75
$ws.Range("J133").Value = 8977.272000000001
$ws.Range("K133").Value = 756023.25
$ws.Range("L133").Value = 26931.816
$ws.Range("M133").Value = -750963.25
$ws.Range("N133").Value = -37051.81600000001
$ws.Range("H136").Value = 14323.625
$ws.Range("I136").Value = 14941.429
$ws.Range("J136").Value = 9999
$ws.Range("K136").Value = 44824.287
$ws.Range("L136").Value = 29997
$ws.Range("M136").Value = -39724.287
$ws.Range("N136").Value = -40197
$ws.Range("H139").Value = 38884.633
$ws.Range("I139").Value = 85567.62
$ws.Range("J139").Value = 3185.8823
$ws.Range("K139").Value = 256702.86
$ws.Range("L139").Value = 9557.6469
$ws.Range("M139").Value = -251562.86
$ws.Range("N139").Value = -19837.6469
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6136.154
$ws.Range("I70").Value = 6167.8066
$ws.Range("J70").Value = 6013.5
$ws.Range("K70").Value = 6167.8066
$ws.Range("L70").Value = 6013.5
$ws.Range("M70").Value = -5897.8066
$ws.Range("N70").Value = -6553.5
$ws.Range("H73").Value = 6136.154
$ws.Range("I73").Value = 6167.8066
$ws.Range("J73").Value = 6013.5
$ws.Range("K73").Value = 6167.8066
$ws.Range("L73").Value = 6013.5
$ws.Range("M73").Value = -5231.8066
$ws.Range("N73").Value = -7885.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8916.916999999999
$ws.Range("I62").Value = 3175
$ws.Range("J62").Value = 11787.875
$ws.Range("K62").Value = 3175
$ws.Range("L62").Value = 11787.875
$ws.Range("M62").Value = -2551
$ws.Range("N62").Value = -13035.875
$ws.Range("H65").Value = 8916.916999999999
$ws.Range("I65").Value = 3175
$ws.Range("J65").Value = 11787.875
$ws.Range("K65").Value = 15875
$ws.Range("L65").Value = 58939.375
$ws.Range("M65").Value = -12755
$ws.Range("N65").Value = -65179.375
